# Technology List.xlsx - add column H with CONCAT formulas building INSERT
# statements from columns C:G, and update the sheet view / selection to
# reflect the state after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column H formulas -----------------------------------------------------
# Row 2 is entered on its own (matches the source which has it as a
# standalone formula, not part of a shared-formula group).
$ws.Range("H2").Formula = "=CONCAT(C2,D2,E2,F2,G2)"

# The remaining rows are filled in three chunks that line up with the
# existing shared-formula breaks already present in column D
# (D3:D66, D67:D130, D131:D165), so that column H gets the same grouping.
$ws.Range("H3:H66").Formula   = "=CONCAT(C3,D3,E3,F3,G3)"
$ws.Range("H67:H130").Formula = "=CONCAT(C67,D67,E67,F67,G67)"
$ws.Range("H131:H165").Formula = "=CONCAT(C131,D131,E131,F131,G131)"

# --- Column H width ----------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 122.43

# --- Sheet view / selection ---------------------------------------------
$app = $ws.Application
$win = $app.ActiveWindow
$win.ScrollRow = 142
$win.ScrollColumn = 3

$ws.Range("H2:H165").Select() | Out-Null
